$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7132.378
$ws.Range("J69").Value = 6999.023
$ws.Range("L69").Value = 20997.069
$ws.Range("N69").Value = -22745.069
$ws.Range("H72").Value = 7132.378
$ws.Range("J72").Value = 6999.023
$ws.Range("L72").Value = 62991.207
$ws.Range("N72").Value = -71727.20699999999
$ws.Range("H80").Value = 942.125
$ws.Range("I80").Value = 496.42856
$ws.Range("K80").Value = 1489.28568
$ws.Range("M80").Value = -491.28568
$ws.Range("H83").Value = 942.125
$ws.Range("I83").Value = 496.42856
$ws.Range("K83").Value = 4467.85704
$ws.Range("M83").Value = 524.1429600000001
$ws.Range("H94").Value = 7401.75
$ws.Range("I94").Value = 7401.75
$ws.Range("K94").Value = 7401.75
$ws.Range("M94").Value = -6950.75
$ws.Range("H107").Value = 107.6
$ws.Range("I107").Value = 79.64286
$ws.Range("K107").Value = 79.64286
$ws.Range("M107").Value = 1840.35714
$ws.Range("H111").Value = 977.7857
$ws.Range("I111").Value = 945.63635
$ws.Range("J111").Value = 1095.6666
$ws.Range("K111").Value = 2836.90905
$ws.Range("L111").Value = 3286.9998
$ws.Range("M111").Value = 230.0909499999998
$ws.Range("N111").Value = -9420.9998
$ws.Range("H115").Value = 1237.4166
$ws.Range("I115").Value = 885
$ws.Range("K115").Value = 2655
$ws.Range("M115").Value = -1088
$ws.Range("H125").Value = 5000
$ws.Range("J125").Value = 5000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -49920

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1668.3846
$ws.Range("I2").Value = 1436.2727
$ws.Range("J2").Value = 2945
$ws.Range("K2").Value = 1436.2727
$ws.Range("L2").Value = 2945
$ws.Range("M2").Value = -1323.2727
$ws.Range("N2").Value = -3171
$ws.Range("H45").Value = 3123.1538
$ws.Range("I45").Value = 2178.2222
$ws.Range("K45").Value = 2178.2222
$ws.Range("M45").Value = -1801.2222
$ws.Range("H61").Value = 2814.625
$ws.Range("I61").Value = 2814.625
$ws.Range("K61").Value = 2814.625
$ws.Range("M61").Value = -2602.625
$ws.Range("H74").Value = 5855.143
$ws.Range("I74").Value = 4735.727
$ws.Range("K74").Value = 4735.727
$ws.Range("M74").Value = -3861.727
$ws.Range("H77").Value = 5855.143
$ws.Range("I77").Value = 4735.727
$ws.Range("K77").Value = 23678.635
$ws.Range("M77").Value = -19310.635
$ws.Range("H116").Value = 1668.3846
$ws.Range("I116").Value = 1436.2727
$ws.Range("J116").Value = 2945
$ws.Range("K116").Value = 1436.2727
$ws.Range("L116").Value = 2945
$ws.Range("M116").Value = 857.7273
$ws.Range("N116").Value = -7533
$ws.Range("H132").Value = 3227.3
$ws.Range("I132").Value = 1063.1666
$ws.Range("K132").Value = 3189.4998
$ws.Range("M132").Value = -659.4998000000001
$ws.Range("H136").Value = 2814.625
$ws.Range("I136").Value = 2814.625
$ws.Range("K136").Value = 8443.875
$ws.Range("M136").Value = -5893.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1668.3846
$ws.Range("I3").Value = 1436.2727
$ws.Range("J3").Value = 2945
$ws.Range("K3").Value = 1436.2727
$ws.Range("L3").Value = 2945
$ws.Range("M3").Value = -1322.2727
$ws.Range("N3").Value = -3173
$ws.Range("H20").Value = 7111.222
$ws.Range("I20").Value = 8024.75
$ws.Range("J20").Value = 5284.1665
$ws.Range("K20").Value = 8024.75
$ws.Range("L20").Value = 5284.1665
$ws.Range("M20").Value = -7777.75
$ws.Range("N20").Value = -5778.1665
$ws.Range("H86").Value = 6394.1816
$ws.Range("I86").Value = 4333.7144
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 4333.7144
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -3210.7144
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 6394.1816
$ws.Range("I89").Value = 4333.7144
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 21668.572
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -16052.572
$ws.Range("N89").Value = -61232
$ws.Range("H107").Value = 3515.1853
$ws.Range("J107").Value = 6727
$ws.Range("L107").Value = 6727
$ws.Range("N107").Value = -10567
$ws.Range("H134").Value = 4047.8333
$ws.Range("I134").Value = 3506.7273
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 10520.1819
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -7985.1819
$ws.Range("N134").Value = -35070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 1863.3
$ws.Range("I86").Value = 1848.1111
$ws.Range("K86").Value = 1848.1111
$ws.Range("M86").Value = -725.1111000000001
$ws.Range("H89").Value = 1863.3
$ws.Range("I89").Value = 1848.1111
$ws.Range("K89").Value = 9240.5555
$ws.Range("M89").Value = -3624.5555
$ws.Range("H132").Value = 3502.4285
$ws.Range("I132").Value = 2420.3333
$ws.Range("K132").Value = 7260.999899999999
$ws.Range("M132").Value = -4730.999899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 121.2
$ws.Range("J12").Value = 123.2
$ws.Range("L12").Value = 369.6
$ws.Range("N12").Value = -715.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 45522.3
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 45522.3
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 45522.3
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -47162.3
$ws.Range("H122").Value = 188208.97
$ws.Range("I122").Value = 280463.06
$ws.Range("J122").Value = 3700.7778
$ws.Range("K122").Value = 841389.1799999999
$ws.Range("L122").Value = 11102.3334
$ws.Range("M122").Value = -838939.1799999999
$ws.Range("N122").Value = -16002.3334
$ws.Range("H132").Value = 46098
$ws.Range("I132").Value = 81746.5
$ws.Range("J132").Value = 7707.3076
$ws.Range("K132").Value = 245239.5
$ws.Range("L132").Value = 23121.9228
$ws.Range("M132").Value = -242709.5
$ws.Range("N132").Value = -28181.9228

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -830
$ws.Range("H20").Value = 9062
$ws.Range("J20").Value = 11999.667
$ws.Range("L20").Value = 11999.667
$ws.Range("N20").Value = -12451.667
$ws.Range("H22").Value = 2138.75
$ws.Range("J22").Value = 3500
$ws.Range("L22").Value = 3500
$ws.Range("N22").Value = -4090
$ws.Range("H27").Value = 2138.75
$ws.Range("J27").Value = 3500
$ws.Range("L27").Value = 3500
$ws.Range("N27").Value = -3714
$ws.Range("H68").Value = 9840.6
$ws.Range("J68").Value = 10050.25
$ws.Range("L68").Value = 10050.25
$ws.Range("N68").Value = -11548.25
$ws.Range("H71").Value = 9840.6
$ws.Range("J71").Value = 10050.25
$ws.Range("L71").Value = 50251.25
$ws.Range("N71").Value = -57739.25
$ws.Range("H82").Value = 4124.077
$ws.Range("I82").Value = 2699.5
$ws.Range("K82").Value = 2699.5
$ws.Range("M82").Value = -2338.5
$ws.Range("H85").Value = 4124.077
$ws.Range("I85").Value = 2699.5
$ws.Range("K85").Value = 2699.5
$ws.Range("M85").Value = -1451.5
$ws.Range("H93").Value = 1526.2667
$ws.Range("I93").Value = 1649.5
$ws.Range("J93").Value = 1033.3334
$ws.Range("K93").Value = 1649.5
$ws.Range("L93").Value = 1033.3334
$ws.Range("M93").Value = -401.5
$ws.Range("N93").Value = -3529.3334
$ws.Range("H132").Value = 4534.1665
$ws.Range("I132").Value = 3928.5
$ws.Range("J132").Value = 5745.5
$ws.Range("K132").Value = 11785.5
$ws.Range("L132").Value = 17236.5
$ws.Range("M132").Value = -9255.5
$ws.Range("N132").Value = -22296.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9142
$ws.Range("J62").Value = 10398.8
$ws.Range("L62").Value = 10398.8
$ws.Range("N62").Value = -11646.8
$ws.Range("H65").Value = 9142
$ws.Range("J65").Value = 10398.8
$ws.Range("L65").Value = 51994
$ws.Range("N65").Value = -58234
$ws.Range("H96").Value = 1000.5
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1001
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1001
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -3747
$ws.Range("H126").Value = 3104.5833
$ws.Range("I126").Value = 1521.875
$ws.Range("K126").Value = 4565.625
$ws.Range("M126").Value = -2095.625
$ws.Range("H136").Value = 2828.9412
$ws.Range("I136").Value = 1719.5186
$ws.Range("J136").Value = 7108.143
$ws.Range("K136").Value = 5158.5558
$ws.Range("L136").Value = 21324.429
$ws.Range("M136").Value = -2608.5558
$ws.Range("N136").Value = -26424.429
